# Applies "Atualizacao de bases das ligas" edit: reorders six already-existing
# Venezuela Primera Division match rows (their B:AC data was shuffled between
# rows while the row index in column A stayed put) and appends four brand new
# match rows (199-202) with the same column layout as the rest of the sheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Step 1: re-shuffle the odds data between existing rows 94-99, 116-117, 135-136 ---
# (column A "id" sequence and columns C/D/E stay as they were on each row; only
#  B and F..AC move to reflect the re-ordered match list)

# Row 94
$ws.Cells.Item(94,2).Value = 6236253
$ws.Cells.Item(94,6).Value = "Deportivo La Guaira"
$ws.Cells.Item(94,7).Value = "UCV"
$ws.Cells.Item(94,8).Value = 0
$ws.Cells.Item(94,9).Value = 0
$ws.Cells.Item(94,10).Value = "D"
$ws.Cells.Item(94,11).Value = 1.833
$ws.Cells.Item(94,12).Value = 3.25
$ws.Cells.Item(94,13).Value = 4
$ws.Cells.Item(94,14).Value = 2
$ws.Cells.Item(94,15).Value = 3.2
$ws.Cells.Item(94,16).Value = 3.5
$ws.Cells.Item(94,17).Value = -0.25
$ws.Cells.Item(94,18).Value = 1.775
$ws.Cells.Item(94,19).Value = 2.025
$ws.Cells.Item(94,20).Value = 2.25
$ws.Cells.Item(94,21).Value = 1.9
$ws.Cells.Item(94,22).Value = 1.9
$ws.Cells.Item(94,23).Value = -1
$ws.Cells.Item(94,24).Value = 2.2
$ws.Cells.Item(94,25).Value = -1
$ws.Cells.Item(94,26).Value = -0.5
$ws.Cells.Item(94,27).Value = 0.5125
$ws.Cells.Item(94,28).Value = -1
$ws.Cells.Item(94,29).Value = 0.8999999999999999

# Row 95
$ws.Cells.Item(95,2).Value = 6236254
$ws.Cells.Item(95,6).Value = "Academia Puerto Cabello"
$ws.Cells.Item(95,7).Value = "Estudiantes Merida"
$ws.Cells.Item(95,8).Value = 1
$ws.Cells.Item(95,9).Value = 0
$ws.Cells.Item(95,10).Value = "H"
$ws.Cells.Item(95,11).Value = 1.727
$ws.Cells.Item(95,12).Value = 3.4
$ws.Cells.Item(95,13).Value = 4.333
$ws.Cells.Item(95,14).Value = 1.666
$ws.Cells.Item(95,15).Value = 3.4
$ws.Cells.Item(95,16).Value = 4.75
$ws.Cells.Item(95,17).Value = -0.75
$ws.Cells.Item(95,18).Value = 1.875
$ws.Cells.Item(95,19).Value = 1.925
$ws.Cells.Item(95,20).Value = 2.5
$ws.Cells.Item(95,21).Value = 1.9
$ws.Cells.Item(95,22).Value = 1.9
$ws.Cells.Item(95,23).Value = 0.6659999999999999
$ws.Cells.Item(95,24).Value = -1
$ws.Cells.Item(95,25).Value = -1
$ws.Cells.Item(95,26).Value = 0.4375
$ws.Cells.Item(95,27).Value = -0.5
$ws.Cells.Item(95,28).Value = -1
$ws.Cells.Item(95,29).Value = 0.8999999999999999

# Row 96
$ws.Cells.Item(96,2).Value = 6236611
$ws.Cells.Item(96,6).Value = "Mineros"
$ws.Cells.Item(96,7).Value = "Monagas"
$ws.Cells.Item(96,8).Value = 2
$ws.Cells.Item(96,9).Value = 1
$ws.Cells.Item(96,10).Value = "H"
$ws.Cells.Item(96,11).Value = 3.2
$ws.Cells.Item(96,12).Value = 3.4
$ws.Cells.Item(96,13).Value = 2
$ws.Cells.Item(96,14).Value = 4.2
$ws.Cells.Item(96,15).Value = 3.8
$ws.Cells.Item(96,16).Value = 1.65
$ws.Cells.Item(96,17).Value = 0.75
$ws.Cells.Item(96,18).Value = 1.95
$ws.Cells.Item(96,19).Value = 1.85
$ws.Cells.Item(96,20).Value = 2.5
$ws.Cells.Item(96,21).Value = 1.825
$ws.Cells.Item(96,22).Value = 1.975
$ws.Cells.Item(96,23).Value = 3.2
$ws.Cells.Item(96,24).Value = -1
$ws.Cells.Item(96,25).Value = -1
$ws.Cells.Item(96,26).Value = 0.95
$ws.Cells.Item(96,27).Value = -1
$ws.Cells.Item(96,28).Value = 0.825
$ws.Cells.Item(96,29).Value = -1

# Row 97
$ws.Cells.Item(97,2).Value = 6236612
$ws.Cells.Item(97,6).Value = "Zamora"
$ws.Cells.Item(97,7).Value = "Carabobo"
$ws.Cells.Item(97,8).Value = 0
$ws.Cells.Item(97,9).Value = 2
$ws.Cells.Item(97,10).Value = "A"
$ws.Cells.Item(97,11).Value = 3.2
$ws.Cells.Item(97,12).Value = 3.1
$ws.Cells.Item(97,13).Value = 2.15
$ws.Cells.Item(97,14).Value = 4.5
$ws.Cells.Item(97,15).Value = 3.3
$ws.Cells.Item(97,16).Value = 1.75
$ws.Cells.Item(97,17).Value = 0.5
$ws.Cells.Item(97,18).Value = 2
$ws.Cells.Item(97,19).Value = 1.8
$ws.Cells.Item(97,20).Value = 2.25
$ws.Cells.Item(97,21).Value = 1.925
$ws.Cells.Item(97,22).Value = 1.875
$ws.Cells.Item(97,23).Value = -1
$ws.Cells.Item(97,24).Value = -1
$ws.Cells.Item(97,25).Value = 0.75
$ws.Cells.Item(97,26).Value = -1
$ws.Cells.Item(97,27).Value = 0.8
$ws.Cells.Item(97,28).Value = -0.5
$ws.Cells.Item(97,29).Value = 0.4375

# Row 98
$ws.Cells.Item(98,2).Value = 6236255
$ws.Cells.Item(98,6).Value = "Deportivo Rayo Zuliano"
$ws.Cells.Item(98,7).Value = "Caracas"
$ws.Cells.Item(98,8).Value = 0
$ws.Cells.Item(98,9).Value = 0
$ws.Cells.Item(98,10).Value = "D"
$ws.Cells.Item(98,11).Value = 3.75
$ws.Cells.Item(98,12).Value = 3.1
$ws.Cells.Item(98,13).Value = 1.95
$ws.Cells.Item(98,14).Value = 2.9
$ws.Cells.Item(98,15).Value = 2.875
$ws.Cells.Item(98,16).Value = 2.45
$ws.Cells.Item(98,17).Value = 0.25
$ws.Cells.Item(98,18).Value = 1.775
$ws.Cells.Item(98,19).Value = 2.025
$ws.Cells.Item(98,20).Value = 2.25
$ws.Cells.Item(98,21).Value = 1.85
$ws.Cells.Item(98,22).Value = 1.95
$ws.Cells.Item(98,23).Value = -1
$ws.Cells.Item(98,24).Value = 1.875
$ws.Cells.Item(98,25).Value = -1
$ws.Cells.Item(98,26).Value = 0.3875
$ws.Cells.Item(98,27).Value = -0.5
$ws.Cells.Item(98,28).Value = -1
$ws.Cells.Item(98,29).Value = 0.95

# Row 99
$ws.Cells.Item(99,2).Value = 6236252
$ws.Cells.Item(99,6).Value = "Deportivo Tachira"
$ws.Cells.Item(99,7).Value = "CD Hermanos Colmenares"
$ws.Cells.Item(99,8).Value = 1
$ws.Cells.Item(99,9).Value = 0
$ws.Cells.Item(99,10).Value = "H"
$ws.Cells.Item(99,11).Value = 1.363
$ws.Cells.Item(99,12).Value = 4.2
$ws.Cells.Item(99,13).Value = 7.5
$ws.Cells.Item(99,14).Value = 1.333
$ws.Cells.Item(99,15).Value = 4.5
$ws.Cells.Item(99,16).Value = 8
$ws.Cells.Item(99,17).Value = -1.5
$ws.Cells.Item(99,18).Value = 2
$ws.Cells.Item(99,19).Value = 1.8
$ws.Cells.Item(99,20).Value = 2.5
$ws.Cells.Item(99,21).Value = 1.925
$ws.Cells.Item(99,22).Value = 1.875
$ws.Cells.Item(99,23).Value = 0.333
$ws.Cells.Item(99,24).Value = -1
$ws.Cells.Item(99,25).Value = -1
$ws.Cells.Item(99,26).Value = -1
$ws.Cells.Item(99,27).Value = 0.8
$ws.Cells.Item(99,28).Value = -1
$ws.Cells.Item(99,29).Value = 0.875

# Row 116
$ws.Cells.Item(116,2).Value = 7352252
$ws.Cells.Item(116,6).Value = "Deportivo Tachira"
$ws.Cells.Item(116,7).Value = "Caracas"
$ws.Cells.Item(116,8).Value = 1
$ws.Cells.Item(116,9).Value = 1
$ws.Cells.Item(116,10).Value = "D"
$ws.Cells.Item(116,11).Value = 2.3
$ws.Cells.Item(116,12).Value = 2.875
$ws.Cells.Item(116,13).Value = 3.1
$ws.Cells.Item(116,14).Value = 2.25
$ws.Cells.Item(116,15).Value = 2.8
$ws.Cells.Item(116,16).Value = 3.25
$ws.Cells.Item(116,17).Value = -0.25
$ws.Cells.Item(116,18).Value = 1.975
$ws.Cells.Item(116,19).Value = 1.825
$ws.Cells.Item(116,20).Value = 2
$ws.Cells.Item(116,21).Value = 1.925
$ws.Cells.Item(116,22).Value = 1.875
$ws.Cells.Item(116,23).Value = -1
$ws.Cells.Item(116,24).Value = 1.8
$ws.Cells.Item(116,25).Value = -1
$ws.Cells.Item(116,26).Value = -0.5
$ws.Cells.Item(116,27).Value = 0.4125
$ws.Cells.Item(116,28).Value = 0
$ws.Cells.Item(116,29).Value = 0

# Row 117
$ws.Cells.Item(117,2).Value = 7352254
$ws.Cells.Item(117,6).Value = "Academia Puerto Cabello"
$ws.Cells.Item(117,7).Value = "Portuguesa"
$ws.Cells.Item(117,8).Value = 1
$ws.Cells.Item(117,9).Value = 1
$ws.Cells.Item(117,10).Value = "D"
$ws.Cells.Item(117,11).Value = 2.05
$ws.Cells.Item(117,12).Value = 3.4
$ws.Cells.Item(117,13).Value = 3
$ws.Cells.Item(117,14).Value = 1.833
$ws.Cells.Item(117,15).Value = 3.5
$ws.Cells.Item(117,16).Value = 3.5
$ws.Cells.Item(117,17).Value = -0.25
$ws.Cells.Item(117,18).Value = 1.65
$ws.Cells.Item(117,19).Value = 2.2
$ws.Cells.Item(117,20).Value = 2.25
$ws.Cells.Item(117,21).Value = 1.825
$ws.Cells.Item(117,22).Value = 1.975
$ws.Cells.Item(117,23).Value = -1
$ws.Cells.Item(117,24).Value = 2.5
$ws.Cells.Item(117,25).Value = -1
$ws.Cells.Item(117,26).Value = -0.5
$ws.Cells.Item(117,27).Value = 0.6000000000000001
$ws.Cells.Item(117,28).Value = -0.5
$ws.Cells.Item(117,29).Value = 0.4875

# Row 135
$ws.Cells.Item(135,2).Value = 7842504
$ws.Cells.Item(135,6).Value = "Angostura FC"
$ws.Cells.Item(135,7).Value = "Deportivo La Guaira"
$ws.Cells.Item(135,8).Value = 1
$ws.Cells.Item(135,9).Value = 1
$ws.Cells.Item(135,10).Value = "D"
$ws.Cells.Item(135,11).Value = 2.75
$ws.Cells.Item(135,12).Value = 3
$ws.Cells.Item(135,13).Value = 2.45
$ws.Cells.Item(135,14).Value = 3.1
$ws.Cells.Item(135,15).Value = 2.875
$ws.Cells.Item(135,16).Value = 2.3
$ws.Cells.Item(135,17).Value = 0.25
$ws.Cells.Item(135,18).Value = 1.8
$ws.Cells.Item(135,19).Value = 2
$ws.Cells.Item(135,20).Value = 2.25
$ws.Cells.Item(135,21).Value = 2.05
$ws.Cells.Item(135,22).Value = 1.75
$ws.Cells.Item(135,23).Value = -1
$ws.Cells.Item(135,24).Value = 1.875
$ws.Cells.Item(135,25).Value = -1
$ws.Cells.Item(135,26).Value = 0.4
$ws.Cells.Item(135,27).Value = -0.5
$ws.Cells.Item(135,28).Value = -0.5
$ws.Cells.Item(135,29).Value = 0.375

# Row 136
$ws.Cells.Item(136,2).Value = 7842507
$ws.Cells.Item(136,6).Value = "Academia Puerto Cabello"
$ws.Cells.Item(136,7).Value = "Estudiantes Merida"
$ws.Cells.Item(136,8).Value = 2
$ws.Cells.Item(136,9).Value = 1
$ws.Cells.Item(136,10).Value = "H"
$ws.Cells.Item(136,11).Value = 1.727
$ws.Cells.Item(136,12).Value = 3.5
$ws.Cells.Item(136,13).Value = 4.2
$ws.Cells.Item(136,14).Value = 1.85
$ws.Cells.Item(136,15).Value = 3.5
$ws.Cells.Item(136,16).Value = 3.6
$ws.Cells.Item(136,17).Value = -0.5
$ws.Cells.Item(136,18).Value = 1.925
$ws.Cells.Item(136,19).Value = 1.875
$ws.Cells.Item(136,20).Value = 2.5
$ws.Cells.Item(136,21).Value = 1.9
$ws.Cells.Item(136,22).Value = 1.9
$ws.Cells.Item(136,23).Value = 0.8500000000000001
$ws.Cells.Item(136,24).Value = -1
$ws.Cells.Item(136,25).Value = -1
$ws.Cells.Item(136,26).Value = 0.925
$ws.Cells.Item(136,27).Value = -1
$ws.Cells.Item(136,28).Value = 0.8999999999999999
$ws.Cells.Item(136,29).Value = -1

# --- Step 2: append four brand-new match rows (199-202) ---

# Row 199
$ws.Cells.Item(199,1).Value = 197
$ws.Cells.Item(199,2).Value = 7977875
$ws.Cells.Item(199,3).Value = "Venezuela Primera Division"
$ws.Cells.Item(199,4).Value = "Venezuela Primera Division"
$ws.Cells.Item(199,5).Value = 45396.66666666666
$ws.Cells.Item(199,6).Value = "Deportivo Rayo Zuliano"
$ws.Cells.Item(199,7).Value = "Deportivo La Guaira"
$ws.Cells.Item(199,8).Value = 2
$ws.Cells.Item(199,9).Value = 2
$ws.Cells.Item(199,10).Value = "D"
$ws.Cells.Item(199,11).Value = 3
$ws.Cells.Item(199,12).Value = 3.1
$ws.Cells.Item(199,13).Value = 2.25
$ws.Cells.Item(199,14).Value = 3.1
$ws.Cells.Item(199,15).Value = 3
$ws.Cells.Item(199,16).Value = 2.3
$ws.Cells.Item(199,17).Value = 0.25
$ws.Cells.Item(199,18).Value = 1.825
$ws.Cells.Item(199,19).Value = 1.975
$ws.Cells.Item(199,20).Value = 2.25
$ws.Cells.Item(199,21).Value = 2
$ws.Cells.Item(199,22).Value = 1.8
$ws.Cells.Item(199,23).Value = -1
$ws.Cells.Item(199,24).Value = 2
$ws.Cells.Item(199,25).Value = -1
$ws.Cells.Item(199,26).Value = 0.4125
$ws.Cells.Item(199,27).Value = -0.5
$ws.Cells.Item(199,28).Value = 1
$ws.Cells.Item(199,29).Value = -1
$ws.Cells.Item(198,1).Copy() | Out-Null
$ws.Cells.Item(199,1).PasteSpecial(-4122) | Out-Null
$ws.Cells.Item(198,5).Copy() | Out-Null
$ws.Cells.Item(199,5).PasteSpecial(-4122) | Out-Null

# Row 200
$ws.Cells.Item(200,1).Value = 198
$ws.Cells.Item(200,2).Value = 7977873
$ws.Cells.Item(200,3).Value = "Venezuela Primera Division"
$ws.Cells.Item(200,4).Value = "Venezuela Primera Division"
$ws.Cells.Item(200,5).Value = 45396.78125
$ws.Cells.Item(200,6).Value = "Deportivo Tachira"
$ws.Cells.Item(200,7).Value = "Carabobo"
$ws.Cells.Item(200,8).Value = 1
$ws.Cells.Item(200,9).Value = 0
$ws.Cells.Item(200,10).Value = "H"
$ws.Cells.Item(200,11).Value = 1.909
$ws.Cells.Item(200,12).Value = 3.2
$ws.Cells.Item(200,13).Value = 3.75
$ws.Cells.Item(200,14).Value = 1.95
$ws.Cells.Item(200,15).Value = 2.9
$ws.Cells.Item(200,16).Value = 3.8
$ws.Cells.Item(200,17).Value = -0.5
$ws.Cells.Item(200,18).Value = 2.025
$ws.Cells.Item(200,19).Value = 1.775
$ws.Cells.Item(200,20).Value = 2
$ws.Cells.Item(200,21).Value = 2
$ws.Cells.Item(200,22).Value = 1.8
$ws.Cells.Item(200,23).Value = 0.95
$ws.Cells.Item(200,24).Value = -1
$ws.Cells.Item(200,25).Value = -1
$ws.Cells.Item(200,26).Value = 1.025
$ws.Cells.Item(200,27).Value = -1
$ws.Cells.Item(200,28).Value = -1
$ws.Cells.Item(200,29).Value = 0.8
$ws.Cells.Item(198,1).Copy() | Out-Null
$ws.Cells.Item(200,1).PasteSpecial(-4122) | Out-Null
$ws.Cells.Item(198,5).Copy() | Out-Null
$ws.Cells.Item(200,5).PasteSpecial(-4122) | Out-Null

# Row 201
$ws.Cells.Item(201,1).Value = 199
$ws.Cells.Item(201,2).Value = 7977872
$ws.Cells.Item(201,3).Value = "Venezuela Primera Division"
$ws.Cells.Item(201,4).Value = "Venezuela Primera Division"
$ws.Cells.Item(201,5).Value = 45396.89583333334
$ws.Cells.Item(201,6).Value = "Caracas"
$ws.Cells.Item(201,7).Value = "UCV"
$ws.Cells.Item(201,8).Value = 2
$ws.Cells.Item(201,9).Value = 2
$ws.Cells.Item(201,10).Value = "D"
$ws.Cells.Item(201,11).Value = 2.45
$ws.Cells.Item(201,12).Value = 2.875
$ws.Cells.Item(201,13).Value = 2.875
$ws.Cells.Item(201,14).Value = 3.3
$ws.Cells.Item(201,15).Value = 2.875
$ws.Cells.Item(201,16).Value = 2.25
$ws.Cells.Item(201,17).Value = 0.25
$ws.Cells.Item(201,18).Value = 1.85
$ws.Cells.Item(201,19).Value = 1.95
$ws.Cells.Item(201,20).Value = 2
$ws.Cells.Item(201,21).Value = 1.975
$ws.Cells.Item(201,22).Value = 1.825
$ws.Cells.Item(201,23).Value = -1
$ws.Cells.Item(201,24).Value = 1.875
$ws.Cells.Item(201,25).Value = -1
$ws.Cells.Item(201,26).Value = 0.425
$ws.Cells.Item(201,27).Value = -0.5
$ws.Cells.Item(201,28).Value = 0.9750000000000001
$ws.Cells.Item(201,29).Value = -1
$ws.Cells.Item(198,1).Copy() | Out-Null
$ws.Cells.Item(201,1).PasteSpecial(-4122) | Out-Null
$ws.Cells.Item(198,5).Copy() | Out-Null
$ws.Cells.Item(201,5).PasteSpecial(-4122) | Out-Null

# Row 202
$ws.Cells.Item(202,1).Value = 200
$ws.Cells.Item(202,2).Value = 7977383
$ws.Cells.Item(202,3).Value = "Venezuela Primera Division"
$ws.Cells.Item(202,4).Value = "Venezuela Primera Division"
$ws.Cells.Item(202,5).Value = 45397.83333333334
$ws.Cells.Item(202,6).Value = "Academia Puerto Cabello"
$ws.Cells.Item(202,7).Value = "Metropolitanos FC"
$ws.Cells.Item(202,8).Value = 2
$ws.Cells.Item(202,9).Value = 3
$ws.Cells.Item(202,10).Value = "A"
$ws.Cells.Item(202,11).Value = 1.666
$ws.Cells.Item(202,12).Value = 3.6
$ws.Cells.Item(202,13).Value = 4.5
$ws.Cells.Item(202,14).Value = 1.7
$ws.Cells.Item(202,15).Value = 3.75
$ws.Cells.Item(202,16).Value = 4.333
$ws.Cells.Item(202,17).Value = -0.75
$ws.Cells.Item(202,18).Value = 1.925
$ws.Cells.Item(202,19).Value = 1.875
$ws.Cells.Item(202,20).Value = 2.25
$ws.Cells.Item(202,21).Value = 1.775
$ws.Cells.Item(202,22).Value = 2.025
$ws.Cells.Item(202,23).Value = -1
$ws.Cells.Item(202,24).Value = -1
$ws.Cells.Item(202,25).Value = 3.333
$ws.Cells.Item(202,26).Value = -1
$ws.Cells.Item(202,27).Value = 0.875
$ws.Cells.Item(202,28).Value = 0.7749999999999999
$ws.Cells.Item(202,29).Value = -1
$ws.Cells.Item(198,1).Copy() | Out-Null
$ws.Cells.Item(202,1).PasteSpecial(-4122) | Out-Null
$ws.Cells.Item(198,5).Copy() | Out-Null
$ws.Cells.Item(202,5).PasteSpecial(-4122) | Out-Null

$excel.CutCopyMode = 0
Write-Host "Applied odds update"